$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Summary paragraph - first bullet: update the two placeholder
#    test-case numbers embedded in the sentence.
# ------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "- Two new test cases (TC-new3333333333333 and TC08888888888888888888-new) were added to validate passwords with less than 8 characters.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Two new test cases (TC-new000000000000000 and TC0111111111111111111-new) were added to validate passwords with less than 8 characters.",
    2)
Write-Host "Step 1 (bullet 1 text) replaced: $found1"

# ------------------------------------------------------------------
# 2. Summary paragraph - second bullet: reword the sentence about the
#    additional section / new requirement.
# ------------------------------------------------------------------
$found2 = $d.Content.Find.Execute(
    "- An additional section (Section 3) was added with new test inputs based on a requirement change from the client to include a special character.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- A new section 3 was added with additional password inputs and expected results based on the new requirement to include a special character.",
    2)
Write-Host "Step 2 (bullet 2 text) replaced: $found2"

# ------------------------------------------------------------------
# 3. Heading further down: "TC" + "-new" + "3333333333333" ->
#    "TC" + "-new" + "000000000000000". The digits live in their own
#    run, sandwiched after the "-new" run. Because this engine
#    coalesces adjacent runs that end up with identical formatting
#    whenever a run's text is edited, we re-stamp the (already
#    matching) resolved font name onto just the replaced span right
#    afterwards; that keeps it as its own distinct run instead of
#    merging into the "-new" run that precedes it.
# ------------------------------------------------------------------
$oldDigits1 = "3333333333333"
$newDigits1 = "000000000000000"
$needle1 = "TC-new" + $oldDigits1
$full = $d.Content.Text
$pos1 = $full.IndexOf($needle1)
if ($pos1 -ge 0) {
    $digitStart1 = $pos1 + "TC-new".Length
    $digitEnd1 = $digitStart1 + $oldDigits1.Length
    $digitRange1 = $d.Range($digitStart1, $digitEnd1)
    $existingFont1 = $digitRange1.Font.Name
    $digitRange1.Text = $newDigits1
    $newDigitRange1 = $d.Range($digitStart1, $digitStart1 + $newDigits1.Length)
    $newDigitRange1.Font.Name = $existingFont1
    Write-Host "Step 3 (TC-new heading digits) replaced at $digitStart1"
} else {
    Write-Host "Step 3 WARNING: needle not found"
}

# ------------------------------------------------------------------
# 4. Heading further down: "TC0" + "8888888888888888888" + "-new" ->
#    "TC0" + "111111111111111111" + "-new". Same technique as step 3,
#    this time the digits run sits between two neighboring runs
#    ("TC0" and "-new") that both share its formatting.
# ------------------------------------------------------------------
$oldDigits2 = "8888888888888888888"
$newDigits2 = "111111111111111111"
$needle2 = "TC0" + $oldDigits2 + "-new"
$full = $d.Content.Text
$pos2 = $full.IndexOf($needle2)
if ($pos2 -ge 0) {
    $digitStart2 = $pos2 + "TC0".Length
    $digitEnd2 = $digitStart2 + $oldDigits2.Length
    $digitRange2 = $d.Range($digitStart2, $digitEnd2)
    $existingFont2 = $digitRange2.Font.Name
    $digitRange2.Text = $newDigits2
    $newDigitRange2 = $d.Range($digitStart2, $digitStart2 + $newDigits2.Length)
    $newDigitRange2.Font.Name = $existingFont2
    Write-Host "Step 4 (TC0...-new heading digits) replaced at $digitStart2"
} else {
    Write-Host "Step 4 WARNING: needle not found"
}
